# DOMA-8315 — fix export of property meters readings.
#
# The "accountNumber" column (column C: {d.i18n.accountNumber} /
# {d.meter[i].accountNumber} / {d.meter[i + 1].accountNumber}) is removed
# from the export template. Deleting the entire column shifts every column
# to its right (D..J) one position to the left (D->C, E->D, ... J->I),
# carrying along cell values/styles and column widths, and shrinks the
# used range from A1:J10 to A1:I10.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").EntireColumn.Delete()
